$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "MCT-2A-Acionamentos Elétricos"

$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "MCT-2A-Acionamentos Elétricos"
$ws.Range("F3").Value = "-"

$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "MCT-3A-Laboratório de Máquinas Elétricas"
$ws.Range("D4").Value = "MCT-3A-Laboratório de Máquinas Elétricas"
$ws.Range("E4").Value = "MCT-2A-Acionamentos Elétricos"
$ws.Range("F4").Value = "-"

$ws.Range("B6").Value = "MCT-2A-Acionamentos Elétricos"
$ws.Range("C6").Value = "MCT-3A-Laboratório de Máquinas Elétricas"
$ws.Range("D6").Value = "MCT-3A-Laboratório de Máquinas Elétricas"
$ws.Range("E6").Value = "MCT-3A-Automação Industrial"
$ws.Range("F6").Value = "-"

$ws.Range("B7").Value = "-"
$ws.Range("D7").Value = "MCT-3A-Laboratório de Máquinas Elétricas"
$ws.Range("E7").Value = "MCT-3A-Automação Industrial"
$ws.Range("F7").Value = "MCT-3A-Automação Industrial"

$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "MCT-3A-Laboratório de Máquinas Elétricas"
$ws.Range("E8").Value = "MCT-3A-Automação Industrial"
